# Apply a reshuffle of the weekly price rows (rows 2-18) in the single
# worksheet. Each data row's entire contents move to a different row
# according to the mapping below (derived from the target diff); row 12
# is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 18
$lastCol = 18   # column R

# Mapping: destination row -> source row (read data from source row in the
# ORIGINAL layout, write it to the destination row in the NEW layout).
$mapping = @{
    2  = 14
    3  = 15
    4  = 8
    5  = 13
    6  = 5
    7  = 6
    8  = 16
    9  = 4
    10 = 2
    11 = 9
    12 = 12
    13 = 17
    14 = 18
    15 = 3
    16 = 11
    17 = 7
    18 = 10
}

# Snapshot all current values before overwriting anything. Use Value2
# (plain, non-parameterized property) since the parameterized Value
# property getter is not reliable for plain reads in this environment.
$snapshot = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $rowVals = @{}
    for ($c = 1; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $rowVals[$c] = $cell.Value2
    }
    $snapshot[$r] = $rowVals
}

# Write the permuted rows back out.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
